$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (RX status) values for the existing data rows 2-10
$ws.Range("C2").Value = "RXA -85.1"
$ws.Range("C3").Value = "RXB -92.8"
$ws.Range("C4").Value = "RXB -94.6"
$ws.Range("C5").Value = "RXA -88"
$ws.Range("C6").Value = "RXA -104"
$ws.Range("C7").Value = "RXA -100.6"
$ws.Range("C8").Value = "RXA -98.1"
$ws.Range("C9").Value = "-"
$ws.Range("C10").Value = "RXA -102.9"

# Remove rows 11-21 which are no longer present in the target data set
$ws.Range("A11:C21").EntireRow.Delete()
